# Insert a new daily price record at row 739 (shifting all subsequent rows
# down by one) for "Pepino ensalada" / Vega Modelo de Temuco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(739).Insert()

$ws.Cells.Item(739, 1).Value = 10
$ws.Cells.Item(739, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(739, 3).Value = "La Araucanía"
$ws.Cells.Item(739, 4).Value = 45124
$ws.Cells.Item(739, 5).Value = 9
$ws.Cells.Item(739, 6).Value = 100112043
$ws.Cells.Item(739, 7).Value = "Pepino ensalada"
$ws.Cells.Item(739, 8).Value = "Sin especificar"
$ws.Cells.Item(739, 9).Value = "Primera"
$ws.Cells.Item(739, 10).Value = 500
$ws.Cells.Item(739, 11).Value = 14000
$ws.Cells.Item(739, 12).Value = 15000
$ws.Cells.Item(739, 13).Value = 14600
$ws.Cells.Item(739, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(739, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(739, 16).Value = 292
$ws.Cells.Item(739, 17).Value = 50
$ws.Cells.Item(739, 18).Value = "Hortaliza"
